$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# temp solve of RWheel: set Fitness column (C) values for rows 2-12 to 4065
$ws.Range("C2:C12").Value = 4065
